$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B27").Value = 9.2252500000000008
$ws.Range("B28").Value = 9.2252500000000008
$ws.Range("B29").Value = 9.2252500000000008
$ws.Range("B30").Value = 9.2252500000000008

$ws.Range("B30").Select()
